$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto listing data (Coin, Link, Price, Volume(1h)) for rows 2-51
$data = @(
    ,@('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '26.456.00', '  +0.17%  ')
    ,@('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.810.68', '  +0.48%  ')
    ,@('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.009', '  +0.03%  ')
    ,@('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.007', '  -0.04%  ')
    ,@('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '306.34', '  -0.62%  ')
    ,@('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4510', '  -0.24%  ')
    ,@('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3589', '  -1.41%  ')
    ,@('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07069', '  -0.09%  ')
    ,@('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.8923', '  +2.84%  ')
    ,@('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07822', '  +0.71%  ')
    ,@('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '19.39', '  +0.79%  ')
    ,@('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.806.48', '  +0.98%  ')
    ,@('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.283', '  +0.69%  ')
    ,@('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.314', '  -0.02%  ')
    ,@('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '85.12', '  -1.02%  ')
    ,@('BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.010', '  -0.06%  ')
    ,@('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000008520', '  -0.23%  ')
    ,@('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.007', '  -0.06%  ')
    ,@('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '26.489.74', '  +0.16%  ')
    ,@('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '14.20', '  +0.15%  ')
    ,@('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.968', '  +0.30%  ')
    ,@('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.050.93', '  +1.65%  ')
    ,@('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '10.52', '  +1.23%  ')
    ,@('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.960', '  -0.68%  ')
    ,@('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '151.75', '  +1.16%  ')
    ,@('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.81', '  -0.27%  ')
    ,@('LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.064', '  +4.12%  ')
    ,@('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '112.30', '  -0.33%  ')
    ,@('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.861', '  +0.30%  ')
    ,@('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08690', '  +0.65%  ')
    ,@('HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.119', '  +3.12%  ')
    ,@('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.807', '  +10.61%  ')
    ,@('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7345', '  +1.04%  ')
    ,@('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.466', '  +0.87%  ')
    ,@('ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.111', '  +0.12%  ')
    ,@('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.077', '  +0.43%  ')
    ,@('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01926', '  +0.67%  ')
    ,@('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05124', '  +1.40%  ')
    ,@('MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.901', '  +1.10%  ')
    ,@('TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5086', '  +3.91%  ')
    ,@('FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.790', '  -2.42%  ')
    ,@('Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1513', '  -3.36%  ')
    ,@('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '8.054', '  -0.56%  ')
    ,@('Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.4674', '  +1.83%  ')
    ,@('PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.008', '  -0.13%  ')
    ,@('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.959', '  +0.66%  ')
    ,@('Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '100.01', '  -0.89%  ')
    ,@('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.574', '  -0.22%  ')
    ,@('Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05991', '  +0.14%  ')
    ,@('Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '63.62', '  +0.46%  ')
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]

    # Price column: values such as "1.010" or "0.4510" must stay as
    # literal text (trailing zeros preserved) instead of being
    # auto-converted to numbers by Excel. Temporarily mark the cell as
    # Text, assign the value, then restore the default (Normal) style
    # so the cell formatting matches the rest of the sheet.
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $entry[2]
    $dCell.Style = "Normal"

    $ws.Cells.Item($row, 5).Value = $entry[3]
    $row++
}
